# Add a "# of cycles" column (L) to the ISA sheet, between the "In presudo verilog"
# column (K) and the "Description" column (old L, now shifted to M). Also refresh
# the saved view state (zoom level / selection) on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "ISA"
$ws2 = $wb.Worksheets.Item(2)   # "presudo instructions"

# --- 1. Insert the new column L on the ISA sheet (shifts old L -> M) ----------
$ws1.Columns("L:L").Insert()

# --- 2. Header cell for the new column -----------------------------------------
$ws1.Range("L1").Value2 = "# of cycles"

# --- 3. Per-instruction cycle counts (rows 2-19) --------------------------------
$cycles = @{
    2  = 1   # add
    3  = 1   # addi
    4  = 1   # sh
    5  = 1   # shi
    6  = 1   # not
    7  = 1   # and
    8  = 1   # or
    9  = 1   # xor
    10 = 1   # cpy
    11 = 1   # cpypc
    12 = 2   # lb
    13 = 2   # sb
    14 = 1   # jmpadr
    15 = 2   # jmpi
    16 = 1   # blt
    17 = 1   # bgt
    18 = 1   # beq
    19 = 1   # bneq
}

foreach ($row in 2..19) {
    $cell = $ws1.Range("L$row")
    $cell.Value2 = $cycles[$row]
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment   = -4108   # xlCenter
    $cell.WrapText = $false
}

# Match the header cell formatting to the rest of the header row.
$ws1.Range("L1").HorizontalAlignment = -4108
$ws1.Range("L1").VerticalAlignment   = -4108
$ws1.Range("L1").WrapText = $true

# --- 4. Blank helper cells alongside the little "for loop" example table -------
$ws1.Range("L22").HorizontalAlignment = -4108
$ws1.Range("L22").VerticalAlignment   = -4108

foreach ($row in 23..26) {
    $cell = $ws1.Range("L$row")
    $cell.HorizontalAlignment = -4131   # xlLeft
    $cell.VerticalAlignment   = -4108   # xlCenter
}

# --- 5. Column widths ------------------------------------------------------------
$ws1.Columns("J").ColumnWidth = 21.29
$ws1.Columns("K").ColumnWidth = 63.43
$ws1.Columns("L").ColumnWidth = 10.43
$ws1.Columns("M").ColumnWidth = 106.57

# --- 6. Refresh the saved view state on sheet 2 first ---------------------------
$ws2.Activate()
$excel.ActiveWindow.Zoom = 70
$ws2.Range("K20").Select()

# --- 7. Refresh the saved view state on sheet 1 (also re-selects it as active) --
$ws1.Activate()
$excel.ActiveWindow.Zoom = 70
$ws1.Range("M23").Select()
